# Add the new data row (row 17) that appears in the diff, appended after
# the existing last row (row 16), and extend the sheet dimension to match.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row = 17

$ws.Cells.Item($row, 1).Value = 11
$ws.Cells.Item($row, 2).Value = 'Vega Monumental Concepción'
$ws.Cells.Item($row, 3).Value = 'Bíobío'
$ws.Cells.Item($row, 4).Value = 45013
$ws.Cells.Item($row, 4).NumberFormat = 'YYYY-MM-DD HH:MM:SS'
$ws.Cells.Item($row, 5).Value = 8
$ws.Cells.Item($row, 6).Value = 'Fruta'
$ws.Cells.Item($row, 7).Value = 100104
$ws.Cells.Item($row, 8).Value = 'Frutos de pepita'
$ws.Cells.Item($row, 9).Value = 100104003
$ws.Cells.Item($row, 10).Value = 'Membrillo'
$ws.Cells.Item($row, 11).Value = 'Champion'
$ws.Cells.Item($row, 12).Value = 'Primera'
$ws.Cells.Item($row, 13).Value = 100
$ws.Cells.Item($row, 14).Value = 9000
$ws.Cells.Item($row, 15).Value = 10000
$ws.Cells.Item($row, 16).Value = 9500
$ws.Cells.Item($row, 17).Value = '$/bandeja 18 kilos granel'
$ws.Cells.Item($row, 18).Value = "Región de O'Higgins"
$ws.Cells.Item($row, 19).Value = 528
$ws.Cells.Item($row, 20).Value = 18
